$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: swap A1 <-> B1 values (id_klienta moves to A1, data moves to B1) ---
$a1 = $ws.Range("A1").Value()
$b1 = $ws.Range("B1").Value()
$ws.Range("A1").Value = $b1
$ws.Range("B1").Value = $a1

# --- Row 2: swap A2 <-> B2 values (0 moves to A2, 2020-12-10 moves to B2) ---
$a2 = $ws.Range("A2").Value()
$b2 = $ws.Range("B2").Value()
$ws.Range("A2").Value = $b2
$ws.Range("B2").Value = $a2

# --- Remove the now-unused A3 / A4 cells entirely ---
$ws.Range("A3").Clear()
$ws.Range("A4").Clear()

# --- Add styled (text-formatted), empty A23 / A24 cells matching column A's style ---
$ws.Range("A23").NumberFormat = "@"
$ws.Range("A24").NumberFormat = "@"

# --- Update the active sheet selection to A1:A2 ---
$ws.Range("A1:A2").Select()
